# Actualización automática 2025-06-09 11:10:09
$wb = $excel.ActiveWorkbook

# ---- Sheet: VENTAS POR GRUPO ----
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsGrupo.Range("D37").Value = 640.34
$wsGrupo.Range("L37").Value = 5640.69
$wsGrupo.Range("D56").Value = "3 de 54"

# ---- Sheet: VENTA MENSUAL ----
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")
$wsMensual.Range("F37").Value = 8023.37
$wsMensual.Range("F56").Value = 22761.77

# ---- Sheet: CUMPLIMIENTO MENSUAL ----
$wsCumplimiento = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
# Column F width: target stored width is 24. This runtime's ColumnWidth
# round-trips with a +5/6 (0.8333...) offset when persisted, so request
# 24 - 0.8333... to land exactly on 24 in the saved file.
$wsCumplimiento.Columns.Item(6).ColumnWidth = 23.166666666666668

$wsCumplimiento.Range("D3").Value = 6284.82
$wsCumplimiento.Range("E3").Value = 7443.18
$wsCumplimiento.Range("F3").Value = 0.4578103146853147

$wsCumplimiento.Range("D16").Value = 7174.84
$wsCumplimiento.Range("E16").Value = 38570.849
$wsCumplimiento.Range("F16").Value = 0.1568418829586325

$wsCumplimiento.Range("D19").Value = 22857.53
$wsCumplimiento.Range("E19").Value = 68105.799
$wsCumplimiento.Range("F19").Value = 0.2512829098416132

$wb.Save()
